$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '69.325.50'
Set-TextValue $ws 'E2' '  +2.56%  '
Set-TextValue $ws 'D3' '3.404.33'
Set-TextValue $ws 'E3' '  +2.58%  '
Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'E4' '  +0.04%  '
Set-TextValue $ws 'D5' '586.52'
Set-TextValue $ws 'E5' '  +1.03%  '
Set-TextValue $ws 'D6' '180.86'
Set-TextValue $ws 'E6' '  +3.90%  '
Set-TextValue $ws 'E7' '  +0.03%  '
Set-TextValue $ws 'E8' '  +1.67%  '
Set-TextValue $ws 'E9' '  +8.97%  '
Set-TextValue $ws 'D10' '0.592'
Set-TextValue $ws 'E10' '  +2.47%  '
Set-TextValue $ws 'D11' '48.56'
Set-TextValue $ws 'E11' '  +4.16%  '
Set-TextValue $ws 'E12' '  +4.59%  '
Set-TextValue $ws 'D13' '685.21'
Set-TextValue $ws 'E13' '  -1.05%  '
Set-TextValue $ws 'D14' '8.66'
Set-TextValue $ws 'D15' '3.958.41'
Set-TextValue $ws 'E15' '  +2.61%  '
Set-TextValue $ws 'D16' '69.503.73'
Set-TextValue $ws 'E16' '  +2.84%  '
Set-TextValue $ws 'D17' '3.410.43'
Set-TextValue $ws 'E17' '  +2.83%  '
Set-TextValue $ws 'E18' '  +1.65%  '
Set-TextValue $ws 'D19' '17.75'
Set-TextValue $ws 'E19' '  +1.64%  '
Set-TextValue $ws 'D20' '11.35'
Set-TextValue $ws 'E20' '  +2.51%  '
Set-TextValue $ws 'D21' '0.911'
Set-TextValue $ws 'E21' '  +2.55%  '
Set-TextValue $ws 'E22' '  +2.69%  '
Set-TextValue $ws 'D23' '5.39'
Set-TextValue $ws 'E23' '  -0.96%  '
Set-TextValue $ws 'D24' '103.07'
Set-TextValue $ws 'E24' '  +1.60%  '
Set-TextValue $ws 'D25' '3.94'
Set-TextValue $ws 'E25' '  +1.38%  '
Set-TextValue $ws 'E26' '  +2.51%  '
Set-TextValue $ws 'D27' '9.71'
Set-TextValue $ws 'E27' '  +3.65%  '
Set-TextValue $ws 'D28' '34.06'
Set-TextValue $ws 'E28' '  +3.85%  '
Set-TextValue $ws 'D29' '8.82'
Set-TextValue $ws 'E29' '  +3.96%  '
Set-TextValue $ws 'D30' '6.96'
Set-TextValue $ws 'E30' '  -0.11%  '
Set-TextValue $ws 'D31' '563.69'
Set-TextValue $ws 'D32' '11.18'
Set-TextValue $ws 'E32' '  +2.01%  '
Set-TextValue $ws 'E33' '  +1.64%  '
Set-TextValue $ws 'D34' '3.58'
Set-TextValue $ws 'E34' '  +10.65%  '
Set-TextValue $ws 'D35' '58.30'
Set-TextValue $ws 'E35' '  +2.00%  '
Set-TextValue $ws 'E36' '  +0.10%  '
Set-TextValue $ws 'D37' '3.670.30'
Set-TextValue $ws 'E37' '  -1.03%  '
Set-TextValue $ws 'D38' '0.141'
Set-TextValue $ws 'E38' '  +6.64%  '
Set-TextValue $ws 'D39' '36.08'
Set-TextValue $ws 'E39' '  +3.38%  '
Set-TextValue $ws 'D40' '0.0₃0721'
Set-TextValue $ws 'E40' '  +8.57%  '
Set-TextValue $ws 'E41' '  +3.89%  '
Set-TextValue $ws 'D42' '2.68'
Set-TextValue $ws 'E42' '  +3.23%  '
Set-TextValue $ws 'D43' '0.340'
Set-TextValue $ws 'E43' '  +2.15%  '
Set-TextValue $ws 'E44' '  +5.40%  '
Set-TextValue $ws 'D45' '3.33'
Set-TextValue $ws 'E45' '  +0.73%  '
Set-TextValue $ws 'D46' '2.69'
Set-TextValue $ws 'E46' '  +2.27%  '
Set-TextValue $ws 'E47' '  +1.54%  '
Set-TextValue $ws 'D48' '1.39'
Set-TextValue $ws 'E48' '  +5.53%  '
Set-TextValue $ws 'E49' '  -0.09%  '
Set-TextValue $ws 'D50' '133.58'
Set-TextValue $ws 'E50' '  +2.00%  '
Set-TextValue $ws 'E51' '  +2.83%  '
